# Update DateBase/orders/Dang Nguyen 195_2026-2-9.xlsx
#
# Appends 10 new order-line rows (51-60) to the "Orders" sheet, pushes the
# old trailing row (A51="1") down to row 61, and refreshes the packed
# "Number" digest string on the "Summary" sheet (G2) so it reflects every
# value now present in the Orders!F column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")
$summary = $wb.Worksheets.Item("Summary")

# Helper: write a value as TEXT (matches the workbook's existing convention
# of storing every cell - including numeric-looking ones - as text with the
# "numberStoredAsText" advisory on). Only needed for cells whose content
# would otherwise auto-detect as a number.
function Set-TextCell {
    param($Sheet, $Address, $Text)
    $cell = $Sheet.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
}

# Helper: write a plain text value (flower names etc. never look numeric,
# so no text-forcing is required - keeps styling untouched, like the rest
# of the sheet).
function Set-Cell {
    param($Sheet, $Address, $Text)
    $Sheet.Range($Address).Value = $Text
}

# ---------------------------------------------------------------------
# 1. Fill in row 51 (previously only had A51 = "1") with real order data,
#    then add the new rows 52-60.
# ---------------------------------------------------------------------
Set-TextCell $ws "A51" "15"
Set-Cell $ws "C51" "147_娜欧米_Red Naomi_Rosa rugosa Thunb._20stems"
Set-TextCell $ws "F51" "17"

Set-Cell $ws "C52" "175_火灵鸟_Free Spirit_Rosa rugosa Thunb._20stems"
Set-TextCell $ws "F52" "4"

Set-Cell $ws "C53" "203_佛罗伊德_Floyd_Rosa rugosa Thunb._20stems"
Set-TextCell $ws "F53" "7.5"

Set-Cell $ws "C54" "411_紫罗兰白_violet white_undefined_1bunch"
Set-TextCell $ws "F54" "20"

Set-TextCell $ws "A55" "16"
Set-Cell $ws "C55" "279_完美甜蜜_undefined_Rosa rugosa Thunb._10stems"
Set-TextCell $ws "F55" "15"

Set-Cell $ws "C56" "144_高原红_High Plateau Red_Rosa rugosa Thunb._20stems"
Set-TextCell $ws "F56" "27"

Set-Cell $ws "C57" "412_紫罗兰粉_violet pink_undefined_1bunch"
Set-TextCell $ws "F57" "10"

Set-TextCell $ws "A58" "17"
Set-Cell $ws "C58" "144_高原红_High Plateau Red_Rosa rugosa Thunb._20stems"
Set-TextCell $ws "F58" "18"

Set-Cell $ws "C59" "203_佛罗伊德_Floyd_Rosa rugosa Thunb._20stems"
Set-TextCell $ws "F59" "18"

Set-Cell $ws "C60" "412_紫罗兰粉_violet pink_undefined_1bunch"
Set-TextCell $ws "F60" "10"

# ---------------------------------------------------------------------
# 2. The row that used to be "51" (just A51 = "1") now lives at row 61.
# ---------------------------------------------------------------------
Set-TextCell $ws "A61" "1"

# ---------------------------------------------------------------------
# 3. Refresh the Summary sheet's packed digest (G2): it is every
#    Orders!F value (F2 downward) concatenated back to back, wrapped with
#    a leading/trailing "0".
# ---------------------------------------------------------------------
Set-TextCell $summary "G2" "03014531467109145105338405302055501059570301001030738510121551542030101530312101051747.5201527101818100"
